$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Unmerge the ranges that are shrinking / disappearing ---
$ws.Range("A3:A6").UnMerge()
$ws.Range("B3:B6").UnMerge()
$ws.Range("A7:A11").UnMerge()
$ws.Range("B7:B11").UnMerge()
$ws.Range("A13:A14").UnMerge()
$ws.Range("B13:B14").UnMerge()

# --- Header cell A1: clear content (was "Nazwa firmy") ---
$ws.Range("A1").Value = ""

# --- Swap "Rodzaj" / "Ilość" header labels on row 2 ---
$ws.Range("C2").Value = "Ilość"
$ws.Range("D2").Value = "Rodzaj"

# --- Row 3 data: new quantity/colour values, "styropian" moves to column D ---
# (write as text formulas then paste-special "values only" so the numeric-looking
#  "2" strings land as real text without picking up a quote-prefix style variant)
$ws.Range("B3").Formula = '="2"'
$ws.Range("B3").Copy()
$ws.Range("B3").PasteSpecial(-4163)
$ws.Range("C3").Formula = '="2"'
$ws.Range("C3").Copy()
$ws.Range("C3").PasteSpecial(-4163)
$ws.Range("D3").Value = "styropian"
$ws.Range("D3").HorizontalAlignment = -4108   # xlCenter - match column B's centred style

# --- Clear out the old colour/size breakdown rows 4-14 entirely ---
$ws.Range("A4:E14").Value = ""
$ws.Range("A4:E14").Style = "Normal"

# --- Re-merge the now single-row groups ---
$ws.Range("A3").Merge()
$ws.Range("B3").Merge()

# --- Column widths: the width-12 custom column moves from C to D ---
$ws.Columns("C").ColumnWidth = 8.43
$ws.Columns("D").ColumnWidth = 11.2
